$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- April (row 41): record Onno's payment and this month's power reading ---
$ws.Range("D41").Value = 147.07
$ws.Range("L41").Value = 100

# Extend the existing shared "down the column" formulas into row 41 / row 42
# (Roy Paid, Roy Balance and Cappy Balance columns hadn't reached these rows yet)
$ws.Range("K41").Formula = "=3 *I41"
$ws.Range("R41").Formula = "=R40 + K41 - (I41 + L41+ M41 + N41 + O41 + P41)"
$ws.Range("W41").Formula = "=W40 + I41 - P41"

# --- May (row 42): the month rolls forward, bringing the formulas with it ---
$ws.Range("E42").Formula = "=D41"
$ws.Range("F42").Value = 59.99
$ws.Range("G42").Formula = "=C42/3"
$ws.Range("H42").Formula = "=(E42+F42)/3"
$ws.Range("I42").Formula = "=G42+H42"
$ws.Range("K42").Formula = "=3 *I42"
$ws.Range("R42").Formula = "=R41 + K42 - (I42 + L42+ M42 + N42 + O42 + P42)"
$ws.Range("W42").Formula = "=W41 + I42 - P42"

# June (row 43) no longer carries an Onno-balance running total past May
$ws.Range("S43").Clear()

# Reflect where the sheet was left scrolled/selected
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 3
$ws.Range("L42").Select()
